# The "Trait lists" column header was renamed to "Tags" on the Template sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Select the cell being edited (matches the resulting selection left behind
# after typing the new header into G1) and update its value.
$ws.Range("G1").Select() | Out-Null
$ws.Range("G1").Value = "Tags"
